# CommonTable.xlsx edit:
#  * Rename sheet "ERROR_CLIENT" -> "ERROR_COMMON"
#  * ERROR_COMMON sheet: drop the stray empty row 4, append LOGIN_SYNC_* rows (66-71)
#  * ERROR_SERVER sheet: rename enum/gen references ErrorServerType -> ServerErrorType,
#    and restore the blank separator row 4

$wb = $excel.ActiveWorkbook

# --- Rename ERROR_CLIENT -> ERROR_COMMON -----------------------------------
$wsCommon = $wb.Worksheets.Item("ERROR_CLIENT")
$wsCommon.Name = "ERROR_COMMON"

# --- ERROR_COMMON sheet body changes ---------------------------------------
# Drop the stray blank row 4 (it already has no cell content; clearing it
# marks the row dirty without adding any cell, so the empty row element is
# dropped entirely on save, and the following rows keep their row numbers).
$wsCommon.Rows.Item(4).ClearContents()

# Append new LOGIN_SYNC_* error rows (66-71)
$newRows = @(
    @{ Row = 66; Id = "LOGIN_SYNC_CANCELLED";          Msg = "Sync cancelled" },
    @{ Row = 67; Id = "LOGIN_SYNC_SLOT_LIST_FAILED";   Msg = "Sync failed while listing slots" },
    @{ Row = 68; Id = "LOGIN_SYNC_LOAD_LOCAL_FAILED";  Msg = "Sync failed while loading local data" },
    @{ Row = 69; Id = "LOGIN_SYNC_LOAD_CLOUD_FAILED";  Msg = "Sync failed while loading cloud data" },
    @{ Row = 70; Id = "LOGIN_SYNC_SAVE_LOCAL_FAILED";  Msg = "Sync failed while saving local data" },
    @{ Row = 71; Id = "LOGIN_SYNC_SAVE_CLOUD_FAILED";  Msg = "Sync failed while saving cloud data" }
)

foreach ($r in $newRows) {
    $wsCommon.Cells.Item($r.Row, 1).Value = $r.Id
    $wsCommon.Cells.Item($r.Row, 3).Value = $r.Msg
}

# --- ERROR_SERVER sheet changes ---------------------------------------------
$wsServer = $wb.Worksheets.Item("ERROR_SERVER")

$wsServer.Cells.Item(2, 1).Value = "enum:ServerErrorType"
$wsServer.Cells.Item(3, 1).Value = "pk, gen:ServerErrorType, code"

# Restore the blank separator row 4 between the schema rows and the data rows.
$wsServer.Rows.Item(4).OutlineLevel = 0
